# Updated generic and feature
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 7 holds a sample SignUp record: C7=Email, D7=FacilityName, E7=DomainURL.
# Swap the placeholder email / facility values for the new ones.
$ws.Range("C7").Value = "prasad868927@yopmail.com"
$ws.Range("D7").Value = "Ganesh998"
$ws.Range("E7").Value = "Ganesh998"

# Move the viewport / selection down to row 7, onto the DomainURL cell.
$ws.Range("E7").Select()
